$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Vwf"
$ws.Range("C2").Value = "Tnfrsf11b"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 56.216147
$ws.Range("H2").Value = 168.648441
$ws.Range("I2").Value = 0.9695233148109074
$ws.Range("J2").Value = 0.9695233148109074
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.06861733333333334
$ws.Range("N2").Value = 0.205852
$ws.Range("O2").Value = 0.01654048691795588
$ws.Range("P2").Value = 0.01654048691795588
$ws.Range("Q2").Value = 3.857402097414667
$ws.Range("R2").Value = 34.716618876732
$ws.Range("S2").Value = 0.01603638770528304
$ws.Range("T2").Value = 0.01603638770528304
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Vwf"
$ws.Range("C3").Value = "Tnfrsf11b"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 56.216147
$ws.Range("H3").Value = 168.648441
$ws.Range("I3").Value = 0.9695233148109074
$ws.Range("J3").Value = 0.9695233148109074
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 3.776574666666666
$ws.Range("N3").Value = 11.329724
$ws.Range("O3").Value = 0.9103586635352137
$ws.Range("P3").Value = 0.9103586635352137
$ws.Range("Q3").Value = 212.3044766178093
$ws.Range("R3").Value = 1910.740289560284
$ws.Range("S3").Value = 0.8826139491374878
$ws.Range("T3").Value = 0.8826139491374878
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Vwf"
$ws.Range("C4").Value = "Tnfrsf11b"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 56.216147
$ws.Range("H4").Value = 168.648441
$ws.Range("I4").Value = 0.9695233148109074
$ws.Range("J4").Value = 0.9695233148109074
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.303255
$ws.Range("N4").Value = 0.909765
$ws.Range("O4").Value = 0.07310084954683041
$ws.Range("P4").Value = 0.07310084954683042
$ws.Range("Q4").Value = 17.047827658485
$ws.Range("R4").Value = 153.430448926365
$ws.Range("S4").Value = 0.07087297796813644
$ws.Range("T4").Value = 0.07087297796813645
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Vwf"
$ws.Range("C5").Value = "Tnfrsf11b"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.07465466666666666
$ws.Range("H5").Value = 0.223964
$ws.Range("I5").Value = 0.001287520467967504
$ws.Range("J5").Value = 0.001287520467967504
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.06861733333333334
$ws.Range("N5").Value = 0.205852
$ws.Range("O5").Value = 0.01654048691795588
$ws.Range("P5").Value = 0.01654048691795588
$ws.Range("Q5").Value = 0.005122604147555556
$ws.Range("R5").Value = 0.046103437328
$ws.Range("S5").Value = 0.00002129621545701694
$ws.Range("T5").Value = 0.00002129621545701695
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Vwf"
$ws.Range("C6").Value = "Tnfrsf11b"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.07465466666666666
$ws.Range("H6").Value = 0.223964
$ws.Range("I6").Value = 0.001287520467967504
$ws.Range("J6").Value = 0.001287520467967504
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 3.776574666666666
$ws.Range("N6").Value = 11.329724
$ws.Range("O6").Value = 0.9103586635352137
$ws.Range("P6").Value = 0.9103586635352137
$ws.Range("Q6").Value = 0.2819389228817777
$ws.Range("R6").Value = 2.537450305936
$ws.Range("S6").Value = 0.00117210541249313
$ws.Range("T6").Value = 0.00117210541249313
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Vwf"
$ws.Range("C7").Value = "Tnfrsf11b"
$ws.Range("D7").Value = "MuSCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 0.07465466666666666
$ws.Range("H7").Value = 0.223964
$ws.Range("I7").Value = 0.001287520467967504
$ws.Range("J7").Value = 0.001287520467967504
$ws.Range("K7").Value = 2
$ws.Range("L7").Value = 0.6666666666666666
$ws.Range("M7").Value = 0.303255
$ws.Range("N7").Value = 0.909765
$ws.Range("O7").Value = 0.07310084954683041
$ws.Range("P7").Value = 0.07310084954683042
$ws.Range("Q7").Value = 0.02263940094
$ws.Range("R7").Value = 0.20375460846
$ws.Range("S7").Value = 0.00009411884001735721
$ws.Range("T7").Value = 0.00009411884001735725
$ws.Range("A8").Value = "MuSCs"
$ws.Range("B8").Value = "Vwf"
$ws.Range("C8").Value = "Tnfrsf11b"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 0.4660483333333333
$ws.Range("H8").Value = 1.398145
$ws.Range("I8").Value = 0.008037632408272877
$ws.Range("J8").Value = 0.008037632408272877
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 0.06861733333333334
$ws.Range("N8").Value = 0.205852
$ws.Range("O8").Value = 0.01654048691795588
$ws.Range("P8").Value = 0.01654048691795588
$ws.Range("Q8").Value = 0.03197899383777778
$ws.Range("R8").Value = 0.28781094454
$ws.Range("S8").Value = 0.0001329463537003758
$ws.Range("T8").Value = 0.0001329463537003758
$ws.Range("A9").Value = "MuSCs"
$ws.Range("B9").Value = "Vwf"
$ws.Range("C9").Value = "Tnfrsf11b"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 0.4660483333333333
$ws.Range("H9").Value = 1.398145
$ws.Range("I9").Value = 0.008037632408272877
$ws.Range("J9").Value = 0.008037632408272877
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 3.776574666666666
$ws.Range("N9").Value = 11.329724
$ws.Range("O9").Value = 0.9103586635352137
$ws.Range("P9").Value = 0.9103586635352137
$ws.Range("Q9").Value = 1.760066329108889
$ws.Range("R9").Value = 15.84059696198
$ws.Range("S9").Value = 0.007317128297182617
$ws.Range("T9").Value = 0.007317128297182617
$ws.Range("A10").Value = "MuSCs"
$ws.Range("B10").Value = "Vwf"
$ws.Range("C10").Value = "Tnfrsf11b"
$ws.Range("D10").Value = "MuSCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 0.4660483333333333
$ws.Range("H10").Value = 1.398145
$ws.Range("I10").Value = 0.008037632408272877
$ws.Range("J10").Value = 0.008037632408272877
$ws.Range("K10").Value = 2
$ws.Range("L10").Value = 0.6666666666666666
$ws.Range("M10").Value = 0.303255
$ws.Range("N10").Value = 0.909765
$ws.Range("O10").Value = 0.07310084954683041
$ws.Range("P10").Value = 0.07310084954683042
$ws.Range("Q10").Value = 0.141331487325
$ws.Range("R10").Value = 1.271983385925
$ws.Range("S10").Value = 0.0005875577573898837
$ws.Range("T10").Value = 0.0005875577573898838
$ws.Range("A11").Value = "Resolving-Mac"
$ws.Range("B11").Value = "Vwf"
$ws.Range("C11").Value = "Tnfrsf11b"
$ws.Range("D11").Value = "ECs"
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 1.226435333333333
$ws.Range("H11").Value = 3.679306
$ws.Range("I11").Value = 0.02115153231285227
$ws.Range("J11").Value = 0.02115153231285228
$ws.Range("K11").Value = 1
$ws.Range("L11").Value = 0.3333333333333333
$ws.Range("M11").Value = 0.06861733333333334
$ws.Range("N11").Value = 0.205852
$ws.Range("O11").Value = 0.01654048691795588
$ws.Range("P11").Value = 0.01654048691795588
$ws.Range("Q11").Value = 0.08415472207911111
$ws.Range("R11").Value = 0.757392498712
$ws.Range("S11").Value = 0.0003498566435154542
$ws.Range("T11").Value = 0.0003498566435154543
$ws.Range("A12").Value = "Resolving-Mac"
$ws.Range("B12").Value = "Vwf"
$ws.Range("C12").Value = "Tnfrsf11b"
$ws.Range("D12").Value = "FAPs"
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 1.226435333333333
$ws.Range("H12").Value = 3.679306
$ws.Range("I12").Value = 0.02115153231285227
$ws.Range("J12").Value = 0.02115153231285228
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 3.776574666666666
$ws.Range("N12").Value = 11.329724
$ws.Range("O12").Value = 0.9103586635352137
$ws.Range("P12").Value = 0.9103586635352137
$ws.Range("Q12").Value = 4.631724610171555
$ws.Range("R12").Value = 41.68552149154399
$ws.Range("S12").Value = 0.01925548068805008
$ws.Range("T12").Value = 0.01925548068805008
$ws.Range("A13").Value = "Resolving-Mac"
$ws.Range("B13").Value = "Vwf"
$ws.Range("C13").Value = "Tnfrsf11b"
$ws.Range("D13").Value = "MuSCs"
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 1.226435333333333
$ws.Range("H13").Value = 3.679306
$ws.Range("I13").Value = 0.02115153231285227
$ws.Range("J13").Value = 0.02115153231285228
$ws.Range("K13").Value = 2
$ws.Range("L13").Value = 0.6666666666666666
$ws.Range("M13").Value = 0.303255
$ws.Range("N13").Value = 0.909765
$ws.Range("O13").Value = 0.07310084954683041
$ws.Range("P13").Value = 0.07310084954683042
$ws.Range("Q13").Value = 0.37192264701
$ws.Range("R13").Value = 3.34730382309
$ws.Range("S13").Value = 0.001546194981286736
$ws.Range("T13").Value = 0.001546194981286736
